# Update results with new concentrations
$wb = $excel.ActiveWorkbook

# --- Sheet "Info": A2/B2 updated ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 4455116336460.7
$wsInfo.Range("B2").Value = 2.180999994277954

# --- Sheet "Activados": rows 2-4 updated, rows 5-20 added (A1:B4 -> A1:B20) ---
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 2; $i -le 20; $i++) {
    $wsAct.Cells.Item($i, 1).Value = 1
    $wsAct.Cells.Item($i, 2).Value = ($i - 2) * 20
}

# --- Sheet "Operando": column A (rows 2-366) changes from 4 to 1 ---
$wsOp = $wb.Worksheets.Item("Operando")
for ($i = 2; $i -le 366; $i++) {
    $wsOp.Cells.Item($i, 1).Value = 1
}

# --- Sheet "Contaminantes": B2:C6 updated ---
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Range("B2").Value = 3750888844080
$wsCont.Range("C2").Value = 139.111
$wsCont.Range("B3").Value = 213683994000
$wsCont.Range("C3").Value = 7.924999999999999
$wsCont.Range("B4").Value = 142784049240
$wsCont.Range("C4").Value = 5.2955
$wsCont.Range("B5").Value = 545340.6966168002
$wsCont.Range("C5").Value = 0.00002022531000000001
$wsCont.Range("B6").Value = 347758903800
$wsCont.Range("C6").Value = 12.8975
